$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.872.62"
$ws.Range("E2").Value = "  +3.22%  "
$ws.Range("D3").Value = "3.632.40"
$ws.Range("E3").Value = "  +6.94%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "588.13"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").Value = "181.39"
$ws.Range("E6").Value = "  +0.70%  "
$ws.Range("D7").Value = "3.627.03"
$ws.Range("E7").Value = "  +7.10%  "
$ws.Range("D8").Value = "0.616"
$ws.Range("E8").Value = "  +3.27%  "
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("E10").Value = "  +1.15%  "
$ws.Range("D11").Value = "0.607"
$ws.Range("E11").Value = "  +2.69%  "
$ws.Range("D12").Value = "49.66"
$ws.Range("E12").Value = "  +2.70%  "
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("D14").Value = "681.68"
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("D15").Value = "4.217.57"
$ws.Range("E15").Value = "  +7.10%  "
$ws.Range("D16").Value = "9.02"
$ws.Range("E16").Value = "  +4.07%  "
$ws.Range("D17").Value = "3.647.79"
$ws.Range("E17").Value = "  +7.71%  "
$ws.Range("D18").Value = "71.938.37"
$ws.Range("E18").Value = "  +3.54%  "
$ws.Range("E19").Value = "  +1.85%  "
$ws.Range("D20").Value = "18.29"
$ws.Range("E20").Value = "  +3.30%  "
$ws.Range("D21").Value = "11.62"
$ws.Range("E21").Value = "  +2.52%  "
$ws.Range("D22").Value = "0.940"
$ws.Range("E22").Value = "  +3.04%  "
$ws.Range("D23").Value = "5.98"
$ws.Range("E23").Value = "  +11.82%  "
$ws.Range("D24").Value = "17.78"
$ws.Range("E24").Value = "  +3.08%  "
$ws.Range("D25").Value = "103.01"
$ws.Range("E25").Value = "  +0.61%  "
$ws.Range("E26").Value = "  +2.18%  "
$ws.Range("E27").Value = "  +4.71%  "
$ws.Range("E28").Value = "  +3.21%  "
$ws.Range("E29").Value = "  +3.84%  "
$ws.Range("D30").Value = "9.21"
$ws.Range("E30").Value = "  +4.45%  "
$ws.Range("E31").Value = "  +5.87%  "
$ws.Range("D32").Value = "4.20"
$ws.Range("E32").Value = "  +15.90%  "
$ws.Range("D33").Value = "585.20"
$ws.Range("E33").Value = "  +5.75%  "
$ws.Range("D34").Value = "11.30"
$ws.Range("E34").Value = "  +1.72%  "
$ws.Range("E35").Value = "  +1.76%  "
$ws.Range("D36").Value = "59.45"
$ws.Range("E36").Value = "  +1.85%  "
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").Value = "3.690.56"
$ws.Range("E38").Value = "  +0.63%  "
$ws.Range("E39").Value = "  +1.72%  "
$ws.Range("D40").Value = "35.71"
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("E41").Value = "  +4.70%  "
$ws.Range("E43").Value = "  +9.45%  "
$ws.Range("E44").Value = "  +2.87%  "
$ws.Range("E45").Value = "  +2.11%  "
$ws.Range("E46").Value = "  +2.87%  "
$ws.Range("D47").Value = "2.81"
$ws.Range("E47").Value = "  +5.33%  "
$ws.Range("E48").Value = "  +2.88%  "
$ws.Range("E49").Value = "  +3.95%  "
$ws.Range("E50").Value = "  -0.29%  "
$ws.Range("D51").Value = "131.75"
$ws.Range("E51").Value = "  +1.71%  "
